$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$textCells = @("D5", "D6", "D10", "D11", "D13", "D19", "D20", "D21", "D22", "D24", "D25", "D27", "D33", "D37", "D38", "D39", "D41", "D42", "D44", "D45", "D48", "D51")
foreach ($addr in $textCells) {
    $ws.Range($addr).NumberFormat = "@"
}

$ws.Range("D2").Value = "70.066.84"
$ws.Range("E2").Value = "  -1.66%  "
$ws.Range("D3").Value = "3.780.45"
$ws.Range("E3").Value = "  +2.76%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "618.22"
$ws.Range("E5").Value = "  +3.32%  "
$ws.Range("D6").Value = "177.77"
$ws.Range("E6").Value = "  -3.64%  "
$ws.Range("D7").Value = "3.779.06"
$ws.Range("E7").Value = "  +2.78%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("E9").Value = "  -0.22%  "
$ws.Range("D10").Value = "0.171"
$ws.Range("E10").Value = "  +4.71%  "
$ws.Range("D11").Value = "6.32"
$ws.Range("E11").Value = "  -3.78%  "
$ws.Range("E12").Value = "  -1.77%  "
$ws.Range("D13").Value = "40.97"
$ws.Range("E13").Value = "  +2.46%  "
$ws.Range("E14").Value = "  +3.03%  "
$ws.Range("D15").Value = "4.410.46"
$ws.Range("E15").Value = "  +2.76%  "
$ws.Range("D16").Value = "3.782.39"
$ws.Range("E16").Value = "  +3.07%  "
$ws.Range("D17").Value = "70.129.39"
$ws.Range("E17").Value = "  -1.50%  "
$ws.Range("E18").Value = "  -0.07%  "
$ws.Range("D19").Value = "7.60"
$ws.Range("E19").Value = "  +1.30%  "
$ws.Range("D20").Value = "16.83"
$ws.Range("E20").Value = "  -0.86%  "
$ws.Range("D21").Value = "509.75"
$ws.Range("E21").Value = "  -1.76%  "
$ws.Range("D22").Value = "9.48"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("E23").Value = "  -2.47%  "
$ws.Range("D24").Value = "87.75"
$ws.Range("E24").Value = "  -0.12%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  +3.23%  "
$ws.Range("E26").Value = "  -3.23%  "
$ws.Range("D27").Value = "10.98"
$ws.Range("E27").Value = "  +1.84%  "
$ws.Range("E28").Value = "  +26.63%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("E30").Value = "  -1.70%  "
$ws.Range("E31").Value = "  +3.19%  "
$ws.Range("E32").Value = "  -4.59%  "
$ws.Range("D33").Value = "31.29"
$ws.Range("E33").Value = "  -1.78%  "
$ws.Range("E34").Value = "  -0.98%  "
$ws.Range("E35").Value = "  +0.19%  "
$ws.Range("E36").Value = "  +4.82%  "
$ws.Range("D37").Value = "6.20"
$ws.Range("E37").Value = "  +0.13%  "
$ws.Range("D38").Value = "0.134"
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("D39").Value = "0.332"
$ws.Range("E39").Value = "  -3.29%  "
$ws.Range("E40").Value = "  -0.17%  "
$ws.Range("D41").Value = "51.00"
$ws.Range("E41").Value = "  +0.50%  "
$ws.Range("D42").Value = "44.94"
$ws.Range("E42").Value = "  -3.07%  "
$ws.Range("E43").Value = "  -0.79%  "
$ws.Range("D44").Value = "417.97"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("D45").Value = "2.82"
$ws.Range("E45").Value = "  +2.28%  "
$ws.Range("D46").Value = "3.031.66"
$ws.Range("E46").Value = "  -4.98%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("D48").Value = "27.42"
$ws.Range("E48").Value = "  -2.89%  "
$ws.Range("E49").Value = "  +2.70%  "
$ws.Range("D51").Value = "2.48"
$ws.Range("E51").Value = "  +1.66%  "
